$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.038.95'
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = '  -3.52%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.809.58'
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").Value = '  +1.10%  '

$ws.Range("E4").Value = '  +0.25%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '593.39'
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").Value = '  -4.46%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '172.33'
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").Value = '  -4.98%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.806.97'
$ws.Range("D7").Style = "Normal"

$ws.Range("E7").Value = '  +1.09%  '

$ws.Range("E8").Value = '  +0.03%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.529'
$ws.Range("D9").Style = "Normal"

$ws.Range("E9").Value = '  -1.36%  '

$ws.Range("E10").Value = '  -6.38%  '

$ws.Range("E11").Value = '  -1.16%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.467'
$ws.Range("D12").Style = "Normal"

$ws.Range("E12").Value = '  -3.44%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '38.30'
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").Value = '  -5.43%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000244'
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").Value = '  -5.71%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.447.15'
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").Value = '  +1.23%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.817.17'
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").Value = '  +1.32%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '68.170.87'
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").Value = '  -3.60%  '

$ws.Range("E18").Value = '  -4.65%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.20'
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").Value = '  -5.45%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.03'
$ws.Range("D20").Style = "Normal"

$ws.Range("E20").Value = '  -3.45%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '488.80'
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").Value = '  -3.67%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.28'
$ws.Range("D22").Style = "Normal"

$ws.Range("E22").Value = '  +0.15%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.734'
$ws.Range("D23").Style = "Normal"

$ws.Range("E23").Value = '  +1.17%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '85.73'
$ws.Range("D24").Style = "Normal"

$ws.Range("E24").Value = '  -1.86%  '

$ws.Range("E25").Value = '  -8.43%  '

$ws.Range("E26").Value = '  -1.79%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.24'
$ws.Range("D27").Style = "Normal"

$ws.Range("E27").Value = '  -6.54%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.14'
$ws.Range("D28").Style = "Normal"

$ws.Range("E28").Value = '  -10.13%  '

$ws.Range("E29").Value = '  -0.28%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.93'
$ws.Range("D30").Style = "Normal"

$ws.Range("E30").Value = '  -0.41%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.43'
$ws.Range("D31").Style = "Normal"

$ws.Range("E31").Value = '  -2.69%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '32.87'
$ws.Range("D32").Style = "Normal"

$ws.Range("E32").Value = '  +6.96%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.61'
$ws.Range("D33").Style = "Normal"

$ws.Range("E33").Value = '  -3.35%  '

$ws.Range("E34").Value = '  -4.53%  '

$ws.Range("E35").Value = '  +0.16%  '

$ws.Range("E36").Value = '  -5.33%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.136'
$ws.Range("D37").Style = "Normal"

$ws.Range("E37").Value = '  -2.64%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.78'
$ws.Range("D38").Style = "Normal"

$ws.Range("E38").Value = '  -5.80%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.324'
$ws.Range("D39").Style = "Normal"

$ws.Range("E39").Value = '  -7.89%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '451.14'
$ws.Range("D40").Style = "Normal"

$ws.Range("E40").Value = '  +3.71%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '48.99'
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").Value = '  -2.04%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.00'
$ws.Range("D42").Style = "Normal"

$ws.Range("E42").Value = '  -4.43%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.88'
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").Value = '  -11.01%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.29'
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").Value = '  -4.28%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '41.22'
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").Value = '  -7.06%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.845.74'
$ws.Range("D46").Style = "Normal"

$ws.Range("E46").Value = '  -4.32%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0352'
$ws.Range("D48").Style = "Normal"

$ws.Range("E48").Value = '  -3.53%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '138.11'
$ws.Range("D49").Style = "Normal"

$ws.Range("E49").Value = '  +1.02%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '26.28'
$ws.Range("D50").Style = "Normal"

$ws.Range("E50").Value = '  -3.95%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '23.60'
$ws.Range("D51").Style = "Normal"

$ws.Range("E51").Value = '  +9.10%  '
